$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Demand_randomforest_Demand_scallg"
$ws.Range("B2").Value = 186.7567334714045
$ws.Range("C2").Value = 0.4932016687630996
$ws.Range("D2").Value = 0.6598485071456115
$ws.Range("E2").Value = "[468.7005020890286, 482.7588393258897, 708.1664568972785, 1257.766019788193]"
$ws.Range("H2").Value = "randomforest"
